# Generate Report for handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) values for the second data row
# (the 878e5e59-... file) on both the "zh-cn" and "de-de" report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-11 13:36:30"
$wsZhCn.Range("G3").Value = "2016-01-11 13:37:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-11 13:36:48"
$wsDeDe.Range("G3").Value = "2016-01-11 13:38:13"
